$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 103482.48
$ws.Range("I15").Value = 103482.48
$ws.Range("K15").Value = 310447.44
$ws.Range("M15").Value = -310278.44
$ws.Range("H18").Value = 1540.4286
$ws.Range("I18").Value = 1356.8
$ws.Range("J18").Value = 1999.5
$ws.Range("K18").Value = 1356.8
$ws.Range("L18").Value = 1999.5
$ws.Range("M18").Value = -1072.8
$ws.Range("N18").Value = -2567.5
$ws.Range("H28").Value = 7346.15
$ws.Range("I28").Value = 519
$ws.Range("J28").Value = 15690.444
$ws.Range("K28").Value = 519
$ws.Range("L28").Value = 15690.444
$ws.Range("M28").Value = -34
$ws.Range("N28").Value = -16660.444
$ws.Range("H32").Value = 15107.75
$ws.Range("J32").Value = 13039.637
$ws.Range("L32").Value = 13039.637
$ws.Range("N32").Value = -13691.637
$ws.Range("H40").Value = 3404
$ws.Range("J40").Value = 3407.6
$ws.Range("L40").Value = 3407.6
$ws.Range("N40").Value = -3757.6
$ws.Range("H43").Value = 1606.4286
$ws.Range("I43").Value = 927.7778
$ws.Range("J43").Value = 2828
$ws.Range("K43").Value = 927.7778
$ws.Range("L43").Value = 2828
$ws.Range("M43").Value = -858.7778
$ws.Range("N43").Value = -2966
$ws.Range("H62").Value = 10424267
$ws.Range("I62").Value = 12828336
$ws.Range("K62").Value = 12828336
$ws.Range("M62").Value = -12827712
$ws.Range("H65").Value = 10424267
$ws.Range("I65").Value = 12828336
$ws.Range("K65").Value = 64141680
$ws.Range("M65").Value = -64138560
$ws.Range("H76").Value = 50005572
$ws.Range("I76").Value = 90914090
$ws.Range("J76").Value = 6272.778
$ws.Range("K76").Value = 90914090
$ws.Range("L76").Value = 6272.778
$ws.Range("M76").Value = -90913775
$ws.Range("N76").Value = -6902.778
$ws.Range("H79").Value = 50005572
$ws.Range("I79").Value = 90914090
$ws.Range("J79").Value = 6272.778
$ws.Range("K79").Value = 90914090
$ws.Range("L79").Value = 6272.778
$ws.Range("M79").Value = -90912998
$ws.Range("N79").Value = -8456.778
$ws.Range("H98").Value = 1376.6786
$ws.Range("I98").Value = 1383.3334
$ws.Range("J98").Value = 1197
$ws.Range("K98").Value = 1383.3334
$ws.Range("L98").Value = 1197
$ws.Range("M98").Value = 114.6666
$ws.Range("N98").Value = -4193
$ws.Range("H101").Value = 4586.6665
$ws.Range("I101").Value = 3895.7144
$ws.Range("K101").Value = 11687.1432
$ws.Range("M101").Value = -10065.1432
$ws.Range("H107").Value = 1494.5814
$ws.Range("I107").Value = 1656.9143
$ws.Range("K107").Value = 1656.9143
$ws.Range("M107").Value = 263.0857000000001
$ws.Range("H116").Value = 6175938
$ws.Range("I116").Value = 10103372
$ws.Range("J116").Value = 4256.5713
$ws.Range("K116").Value = 10103372
$ws.Range("L116").Value = 4256.5713
$ws.Range("M116").Value = -10099930
$ws.Range("N116").Value = -11140.5713
$ws.Range("H122").Value = 1376.6786
$ws.Range("I122").Value = 1383.3334
$ws.Range("J122").Value = 1197
$ws.Range("K122").Value = 4150.0002
$ws.Range("L122").Value = 3591
$ws.Range("M122").Value = -1700.0002
$ws.Range("N122").Value = -8491
$ws.Range("H125").Value = 2377.1
$ws.Range("J125").Value = 922.4
$ws.Range("L125").Value = 8301.6
$ws.Range("N125").Value = -13221.6
$ws.Range("H137").Value = 3223.5
$ws.Range("I137").Value = 1739.8334
$ws.Range("K137").Value = 5219.5002
$ws.Range("M137").Value = -2669.5002
$ws.Range("H138").Value = 2824.0408
$ws.Range("I138").Value = 1027.1471
$ws.Range("J138").Value = 3778.6406
$ws.Range("K138").Value = 3081.4413
$ws.Range("L138").Value = 11335.9218
$ws.Range("M138").Value = 2058.5587
$ws.Range("N138").Value = -21615.9218
$ws.Range("H141").Value = 2946
$ws.Range("I141").Value = 2204.818
$ws.Range("J141").Value = 4576.6
$ws.Range("K141").Value = 6614.454000000001
$ws.Range("L141").Value = 13729.8
$ws.Range("M141").Value = -1434.454000000001
$ws.Range("N141").Value = -24089.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3248.2239
$ws.Range("I2").Value = 1478.6346
$ws.Range("J2").Value = 9382.799999999999
$ws.Range("K2").Value = 1478.6346
$ws.Range("L2").Value = 9382.799999999999
$ws.Range("M2").Value = -1365.6346
$ws.Range("N2").Value = -9608.799999999999
$ws.Range("H5").Value = 1313.6666
$ws.Range("I5").Value = 1313.6666
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1313.6666
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1201.6666
$ws.Range("N5").ClearContents()
$ws.Range("H32").Value = 3184022.5
$ws.Range("I32").Value = 8410.207
$ws.Range("K32").Value = 8410.207
$ws.Range("M32").Value = -8123.207
$ws.Range("H45").Value = 2606.5264
$ws.Range("I45").Value = 2518
$ws.Range("J45").Value = 4200
$ws.Range("K45").Value = 2518
$ws.Range("L45").Value = 4200
$ws.Range("M45").Value = -2141
$ws.Range("N45").Value = -4954
$ws.Range("H61").Value = 5073.6055
$ws.Range("I61").Value = 5085.273
$ws.Range("J61").Value = 4996.6
$ws.Range("K61").Value = 5085.273
$ws.Range("L61").Value = 4996.6
$ws.Range("M61").Value = -4873.273
$ws.Range("N61").Value = -5420.6
$ws.Range("H74").Value = 5433.68
$ws.Range("I74").Value = 5330.1333
$ws.Range("K74").Value = 5330.1333
$ws.Range("M74").Value = -4456.1333
$ws.Range("H77").Value = 5433.68
$ws.Range("I77").Value = 5330.1333
$ws.Range("K77").Value = 26650.6665
$ws.Range("M77").Value = -22282.6665
$ws.Range("H110").Value = 5985.4644
$ws.Range("I110").Value = 2986.3333
$ws.Range("J110").Value = 7406.1055
$ws.Range("K110").Value = 2986.3333
$ws.Range("L110").Value = 7406.1055
$ws.Range("M110").Value = -941.3332999999998
$ws.Range("N110").Value = -11496.1055
$ws.Range("H116").Value = 3248.2239
$ws.Range("I116").Value = 1478.6346
$ws.Range("J116").Value = 9382.799999999999
$ws.Range("K116").Value = 1478.6346
$ws.Range("L116").Value = 9382.799999999999
$ws.Range("M116").Value = 815.3653999999999
$ws.Range("N116").Value = -13970.8
$ws.Range("H132").Value = 752953.8
$ws.Range("I132").Value = 871048.8
$ws.Range("J132").Value = 123114
$ws.Range("K132").Value = 2613146.4
$ws.Range("L132").Value = 369342
$ws.Range("M132").Value = -2610616.4
$ws.Range("N132").Value = -374402
$ws.Range("H136").Value = 5073.6055
$ws.Range("I136").Value = 5085.273
$ws.Range("J136").Value = 4996.6
$ws.Range("K136").Value = 15255.819
$ws.Range("L136").Value = 14989.8
$ws.Range("M136").Value = -12705.819
$ws.Range("N136").Value = -20089.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3248.2239
$ws.Range("I3").Value = 1478.6346
$ws.Range("J3").Value = 9382.799999999999
$ws.Range("K3").Value = 1478.6346
$ws.Range("L3").Value = 9382.799999999999
$ws.Range("M3").Value = -1364.6346
$ws.Range("N3").Value = -9610.799999999999
$ws.Range("H4").Value = 1313.6666
$ws.Range("I4").Value = 1313.6666
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 1313.6666
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -1198.6666
$ws.Range("N4").ClearContents()
$ws.Range("H94").Value = 5023.4863
$ws.Range("I94").Value = 3733.1538
$ws.Range("K94").Value = 3733.1538
$ws.Range("M94").Value = -3282.1538
$ws.Range("H99").Value = 7535.423
$ws.Range("I99").Value = 7527.1333
$ws.Range("J99").Value = 7546.727
$ws.Range("K99").Value = 7527.1333
$ws.Range("L99").Value = 7546.727
$ws.Range("M99").Value = -6029.1333
$ws.Range("N99").Value = -10542.727
$ws.Range("H105").Value = 2233.75
$ws.Range("I105").Value = 2298.7368
$ws.Range("K105").Value = 2298.7368
$ws.Range("M105").Value = -551.7368000000001
$ws.Range("H107").Value = 7153514.5
$ws.Range("I107").Value = 11119627
$ws.Range("K107").Value = 11119627
$ws.Range("M107").Value = -11117707
$ws.Range("H134").Value = 749073.5600000001
$ws.Range("I134").Value = 909138
$ws.Range("J134").Value = 5917.2144
$ws.Range("K134").Value = 2727414
$ws.Range("L134").Value = 17751.6432
$ws.Range("M134").Value = -2724879
$ws.Range("N134").Value = -22821.6432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 11630607
$ws.Range("I16").Value = 14287098
$ws.Range("J16").Value = 8460.125
$ws.Range("K16").Value = 14287098
$ws.Range("L16").Value = 8460.125
$ws.Range("M16").Value = -14286811
$ws.Range("N16").Value = -9034.125
$ws.Range("H31").Value = 8156.0493
$ws.Range("I31").Value = 11415.789
$ws.Range("J31").Value = 6681.405
$ws.Range("K31").Value = 11415.789
$ws.Range("L31").Value = 6681.405
$ws.Range("M31").Value = -11120.789
$ws.Range("N31").Value = -7271.405
$ws.Range("H34").Value = 8156.0493
$ws.Range("I34").Value = 11415.789
$ws.Range("J34").Value = 6681.405
$ws.Range("K34").Value = 11415.789
$ws.Range("L34").Value = 6681.405
$ws.Range("M34").Value = -11213.789
$ws.Range("N34").Value = -7085.405
$ws.Range("H50").Value = 30039.588
$ws.Range("I50").Value = 13333.333
$ws.Range("J50").Value = 33619.5
$ws.Range("K50").Value = 13333.333
$ws.Range("L50").Value = 33619.5
$ws.Range("M50").Value = -12708.333
$ws.Range("N50").Value = -34869.5
$ws.Range("H58").Value = 7043.393
$ws.Range("I58").Value = 4030.2727
$ws.Range("J58").Value = 18091.5
$ws.Range("K58").Value = 4030.2727
$ws.Range("L58").Value = 18091.5
$ws.Range("M58").Value = -3827.2727
$ws.Range("N58").Value = -18497.5
$ws.Range("H74").Value = 16877.8
$ws.Range("J74").Value = 16877.8
$ws.Range("L74").Value = 16877.8
$ws.Range("N74").Value = -18625.8
$ws.Range("H77").Value = 16877.8
$ws.Range("J77").Value = 16877.8
$ws.Range("L77").Value = 50633.39999999999
$ws.Range("N77").Value = -59369.39999999999
$ws.Range("H86").Value = 7951.522
$ws.Range("I86").Value = 4770
$ws.Range("J86").Value = 12900.556
$ws.Range("K86").Value = 4770
$ws.Range("L86").Value = 12900.556
$ws.Range("M86").Value = -3647
$ws.Range("N86").Value = -15146.556
$ws.Range("H89").Value = 7951.522
$ws.Range("I89").Value = 4770
$ws.Range("J89").Value = 12900.556
$ws.Range("K89").Value = 23850
$ws.Range("L89").Value = 64502.78
$ws.Range("M89").Value = -18234
$ws.Range("N89").Value = -75734.78
$ws.Range("H93").Value = 20317.666
$ws.Range("I93").Value = 20317.666
$ws.Range("K93").Value = 20317.666
$ws.Range("M93").Value = -18445.666
$ws.Range("H99").Value = 7941355.5
$ws.Range("I99").Value = 13893684
$ws.Range("K99").Value = 13893684
$ws.Range("M99").Value = -13892186
$ws.Range("H103").Value = 10524
$ws.Range("I103").Value = 10524
$ws.Range("K103").Value = 10524
$ws.Range("M103").Value = -9352
$ws.Range("H105").Value = 43479384
$ws.Range("J105").Value = 1193.5
$ws.Range("L105").Value = 1193.5
$ws.Range("N105").Value = -4687.5
$ws.Range("H106").Value = 34368.5
$ws.Range("J106").Value = 35828.332
$ws.Range("L106").Value = 35828.332
$ws.Range("N106").Value = -38352.332
$ws.Range("H113").Value = 11630607
$ws.Range("I113").Value = 14287098
$ws.Range("J113").Value = 8460.125
$ws.Range("K113").Value = 14287098
$ws.Range("L113").Value = 8460.125
$ws.Range("M113").Value = -14284928
$ws.Range("N113").Value = -12800.125
$ws.Range("H122").Value = 8195.75
$ws.Range("I122").Value = 3172.889
$ws.Range("J122").Value = 23264.334
$ws.Range("K122").Value = 9518.667000000001
$ws.Range("L122").Value = 69793.00199999999
$ws.Range("M122").Value = -7068.667000000001
$ws.Range("N122").Value = -74693.00199999999
$ws.Range("H126").Value = 7941355.5
$ws.Range("I126").Value = 13893684
$ws.Range("K126").Value = 41681052
$ws.Range("M126").Value = -41678582
$ws.Range("H132").Value = 5345.439
$ws.Range("I132").Value = 2961.889
$ws.Range("J132").Value = 22507
$ws.Range("K132").Value = 8885.667000000001
$ws.Range("L132").Value = 67521
$ws.Range("M132").Value = -6355.667000000001
$ws.Range("N132").Value = -72581
$ws.Range("H134").Value = 11302.5
$ws.Range("I134").Value = 6199
$ws.Range("J134").Value = 17681.875
$ws.Range("K134").Value = 18597
$ws.Range("L134").Value = 53045.625
$ws.Range("M134").Value = -16062
$ws.Range("N134").Value = -58115.625
$ws.Range("H136").Value = 7043.393
$ws.Range("I136").Value = 4030.2727
$ws.Range("J136").Value = 18091.5
$ws.Range("K136").Value = 12090.8181
$ws.Range("L136").Value = 54274.5
$ws.Range("M136").Value = -9540.8181
$ws.Range("N136").Value = -59374.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 20086
$ws.Range("I3").Value = 3088.2856
$ws.Range("J3").Value = 39916.668
$ws.Range("K3").Value = 9264.856800000001
$ws.Range("L3").Value = 119750.004
$ws.Range("M3").Value = -9152.856800000001
$ws.Range("N3").Value = -119974.004
$ws.Range("H5").Value = 1216.7727
$ws.Range("I5").Value = 633.2857
$ws.Range("J5").Value = 1749.5217
$ws.Range("K5").Value = 1899.8571
$ws.Range("L5").Value = 5248.5651
$ws.Range("M5").Value = -1787.8571
$ws.Range("N5").Value = -5472.5651
$ws.Range("H34").Value = 2839.1333
$ws.Range("I34").Value = 253.5
$ws.Range("J34").Value = 8010.4
$ws.Range("K34").Value = 760.5
$ws.Range("L34").Value = 24031.2
$ws.Range("M34").Value = -676.5
$ws.Range("N34").Value = -24199.2
$ws.Range("H37").Value = 120300
$ws.Range("J37").Value = 120300
$ws.Range("L37").Value = 360900
$ws.Range("N37").Value = -361124
$ws.Range("H55").Value = 796066.9
$ws.Range("I55").Value = 44699.43
$ws.Range("J55").Value = 1672662.1
$ws.Range("K55").Value = 134098.29
$ws.Range("L55").Value = 5017986.300000001
$ws.Range("M55").Value = -133921.29
$ws.Range("N55").Value = -5018340.300000001
$ws.Range("H63").Value = 16977.4
$ws.Range("J63").Value = 19971.75
$ws.Range("L63").Value = 59915.25
$ws.Range("N63").Value = -61413.25
$ws.Range("H66").Value = 16977.4
$ws.Range("J66").Value = 19971.75
$ws.Range("L66").Value = 179745.75
$ws.Range("N66").Value = -187233.75
$ws.Range("H101").Value = 35004.4
$ws.Range("J101").Value = 39999
$ws.Range("L101").Value = 119997
$ws.Range("N101").Value = -124865
$ws.Range("H113").Value = 2977343.5
$ws.Range("J113").Value = 1377
$ws.Range("L113").Value = 4131
$ws.Range("N113").Value = -8471
$ws.Range("H130").Value = 10676.333
$ws.Range("J130").Value = 14999.5
$ws.Range("L130").Value = 44998.5
$ws.Range("N130").Value = -55038.5
$ws.Range("H131").Value = 3387
$ws.Range("I131").Value = 3611.5715
$ws.Range("K131").Value = 10834.7145
$ws.Range("M131").Value = -5794.7145
$ws.Range("H135").Value = 1216.7727
$ws.Range("I135").Value = 633.2857
$ws.Range("J135").Value = 1749.5217
$ws.Range("K135").Value = 5699.571300000001
$ws.Range("L135").Value = 15745.6953
$ws.Range("M135").Value = -3164.571300000001
$ws.Range("N135").Value = -20815.6953

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5040.154
$ws.Range("I70").Value = 4865.636
$ws.Range("J70").Value = 6000
$ws.Range("K70").Value = 4865.636
$ws.Range("L70").Value = 6000
$ws.Range("M70").Value = -4595.636
$ws.Range("N70").Value = -6540
$ws.Range("H73").Value = 5040.154
$ws.Range("I73").Value = 4865.636
$ws.Range("J73").Value = 6000
$ws.Range("K73").Value = 4865.636
$ws.Range("L73").Value = 6000
$ws.Range("M73").Value = -3929.636
$ws.Range("N73").Value = -7872
$ws.Range("H80").Value = 6332.7
$ws.Range("I80").Value = 3516.2
$ws.Range("J80").Value = 9149.200000000001
$ws.Range("K80").Value = 3516.2
$ws.Range("L80").Value = 9149.200000000001
$ws.Range("M80").Value = -2518.2
$ws.Range("N80").Value = -11145.2
$ws.Range("H83").Value = 6332.7
$ws.Range("I83").Value = 3516.2
$ws.Range("J83").Value = 9149.200000000001
$ws.Range("K83").Value = 17581
$ws.Range("L83").Value = 45746
$ws.Range("M83").Value = -12589
$ws.Range("N83").Value = -55730
$ws.Range("H102").Value = 3987.971
$ws.Range("I102").Value = 2534.8572
$ws.Range("J102").Value = 7548.1
$ws.Range("K102").Value = 2534.8572
$ws.Range("L102").Value = 7548.1
$ws.Range("M102").Value = -912.8571999999999
$ws.Range("N102").Value = -10792.1
$ws.Range("H113").Value = 7574.2
$ws.Range("I113").Value = 3808
$ws.Range("J113").Value = 10085
$ws.Range("K113").Value = 3808
$ws.Range("L113").Value = 10085
$ws.Range("M113").Value = -1638
$ws.Range("N113").Value = -14425
$ws.Range("H122").Value = 4486.8667
$ws.Range("I122").Value = 3660.9
$ws.Range("J122").Value = 6138.8
$ws.Range("K122").Value = 10982.7
$ws.Range("L122").Value = 18416.4
$ws.Range("M122").Value = -8532.700000000001
$ws.Range("N122").Value = -23316.4
$ws.Range("H132").Value = 4025.4558
$ws.Range("I132").Value = 3270.653
$ws.Range("K132").Value = 9811.958999999999
$ws.Range("M132").Value = -7281.958999999999
$ws.Range("H136").Value = 20931.2
$ws.Range("J136").Value = 20931.2
$ws.Range("L136").Value = 62793.60000000001
$ws.Range("N136").Value = -67893.60000000001
$ws.Range("H140").Value = 65162.125
$ws.Range("J140").Value = 65162.125
$ws.Range("L140").Value = 65162.125
$ws.Range("N140").Value = -75522.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6603.645
$ws.Range("I7").Value = 6802.136
$ws.Range("J7").Value = 6118.4443
$ws.Range("K7").Value = 6802.136
$ws.Range("L7").Value = 6118.4443
$ws.Range("M7").Value = -6690.136
$ws.Range("N7").Value = -6342.4443
$ws.Range("H22").Value = 667.93335
$ws.Range("I22").Value = 604.44446
$ws.Range("J22").Value = 763.1667
$ws.Range("K22").Value = 604.44446
$ws.Range("L22").Value = 763.1667
$ws.Range("M22").Value = -309.44446
$ws.Range("N22").Value = -1353.1667
$ws.Range("H25").Value = 3336335.8
$ws.Range("I25").Value = 3336335.8
$ws.Range("K25").Value = 3336335.8
$ws.Range("M25").Value = -3336105.8
$ws.Range("H27").Value = 667.93335
$ws.Range("I27").Value = 604.44446
$ws.Range("J27").Value = 763.1667
$ws.Range("K27").Value = 604.44446
$ws.Range("L27").Value = 763.1667
$ws.Range("M27").Value = -497.44446
$ws.Range("N27").Value = -977.1667
$ws.Range("H40").Value = 3041.879
$ws.Range("I40").Value = 2499.2593
$ws.Range("J40").Value = 5483.6665
$ws.Range("K40").Value = 2499.2593
$ws.Range("L40").Value = 5483.6665
$ws.Range("M40").Value = -2363.2593
$ws.Range("N40").Value = -5755.6665
$ws.Range("H46").Value = 21740248
$ws.Range("J46").Value = 41667976
$ws.Range("L46").Value = 41667976
$ws.Range("N46").Value = -41668352
$ws.Range("H47").Value = 30059
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H52").Value = 30059
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H55").Value = 3632.7297
$ws.Range("I55").Value = 2594.889
$ws.Range("J55").Value = 6434.9
$ws.Range("K55").Value = 2594.889
$ws.Range("L55").Value = 6434.9
$ws.Range("M55").Value = -2421.889
$ws.Range("N55").Value = -6780.9
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()
$ws.Range("H122").Value = 458811.38
$ws.Range("I122").Value = 530098.9
$ws.Range("K122").Value = 1590296.7
$ws.Range("M122").Value = -1587846.7
$ws.Range("H126").Value = 6603.645
$ws.Range("I126").Value = 6802.136
$ws.Range("J126").Value = 6118.4443
$ws.Range("K126").Value = 20406.408
$ws.Range("L126").Value = 18355.3329
$ws.Range("M126").Value = -17936.408
$ws.Range("N126").Value = -23295.3329
$ws.Range("H132").Value = 4895.4736
$ws.Range("I132").Value = 3528.825
$ws.Range("K132").Value = 10586.475
$ws.Range("M132").Value = -8056.474999999999
$ws.Range("H136").Value = 6171.892
$ws.Range("I136").Value = 4141.615
$ws.Range("K136").Value = 12424.845
$ws.Range("M136").Value = -9874.844999999999
$ws.Range("H140").Value = 71750
$ws.Range("J140").Value = 74000
$ws.Range("L140").Value = 74000
$ws.Range("N140").Value = -84360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 25010000
$ws.Range("J4").Value = 50000000
$ws.Range("L4").Value = 50000000
$ws.Range("N4").Value = -50000226
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("H81").Value = 1917.3572
$ws.Range("I81").Value = 1607.5
$ws.Range("J81").Value = 2149.75
$ws.Range("K81").Value = 3215
$ws.Range("L81").Value = 4299.5
$ws.Range("M81").Value = -2154
$ws.Range("N81").Value = -6421.5
$ws.Range("H84").Value = 1917.3572
$ws.Range("I84").Value = 1607.5
$ws.Range("J84").Value = 2149.75
$ws.Range("K84").Value = 16075
$ws.Range("L84").Value = 21497.5
$ws.Range("M84").Value = -10771
$ws.Range("N84").Value = -32105.5
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H113").Value = 5051701
$ws.Range("I113").Value = 8773130
$ws.Range("J113").Value = 1189.8572
$ws.Range("K113").Value = 26319390
$ws.Range("L113").Value = 3569.5716
$ws.Range("M113").Value = -26317220
$ws.Range("N113").Value = -7909.571599999999
$ws.Range("H122").Value = 7205.7295
$ws.Range("I122").Value = 2105.08
$ws.Range("K122").Value = 6315.24
$ws.Range("M122").Value = -3865.24
$ws.Range("H126").Value = 5011.5
$ws.Range("I126").Value = 2710.3333
$ws.Range("J126").Value = 11915
$ws.Range("K126").Value = 8130.999899999999
$ws.Range("L126").Value = 35745
$ws.Range("M126").Value = -5660.999899999999
$ws.Range("N126").Value = -40685
$ws.Range("H132").Value = 4832.1343
$ws.Range("I132").Value = 4128.771
$ws.Range("J132").Value = 6609.0527
$ws.Range("K132").Value = 12386.313
$ws.Range("L132").Value = 19827.1581
$ws.Range("M132").Value = -9856.312999999998
$ws.Range("N132").Value = -24887.1581
$ws.Range("H136").Value = 7150109.5
$ws.Range("I136").Value = 11370223
$ws.Range("J136").Value = 8379.885
$ws.Range("K136").Value = 34110669
$ws.Range("L136").Value = 25139.655
$ws.Range("M136").Value = -34108119
$ws.Range("N136").Value = -30239.655
